$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.751.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.664.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.55"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.12"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.658"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +6.51%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.01%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.83"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.81"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000194"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.142.81"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.617.85"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.650.06"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.57"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.79"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "349.39"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.42"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.75"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.80"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +9.31%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.52"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.63"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "565.05"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.12%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.06"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.162"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.13"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.82"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.70"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.44"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.421"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.55"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "155.57"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.92"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "160.61"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0609"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.28"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.70"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.639"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.102"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0254"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.78"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0244"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.92%  "
